$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "image" column (BL) appended after the existing last column (BK).
# Header cell BL1 gets the same bold/center/top formatting as the other
# header cells, plus a thin left/right border (instead of the all-round
# thin border used by the rest of row 1).
$header = $ws.Range("BL1")
$header.Value = "image"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.Item(7).LineStyle = 1
$header.Borders.Item(10).LineStyle = 1

# Data rows 2-5 hold the image file names associated with each IOA row.
$ws.Range("BL2").Value = "image3.png"
$ws.Range("BL3").Value = "image12.png"
$ws.Range("BL4").Value = "image8.png"
$ws.Range("BL5").Value = "image10.png"
